{"js": "// Replace the division-problem answers in the table with the new set\n// of problems, per the commit's regenerated output.\nconst replacements = [\n  [\"509\u00f74=127, 1\", \"152\u00f73=50, 2\"],\n  [\"961\u00f79=106, 7\", \"884\u00f77=126, 2\"],\n  [\"364\u00f79=40, 4\", \"163\u00f76=27, 1\"],\n  [\"693\u00f78=86, 5\", \"546\u00f79=60, 6\"],\n  [\"725\u00f72=362, 1\", \"916\u00f75=183, 1\"],\n  [\"766\u00f72=383, 0\", \"875\u00f79=97, 2\"],\n  [\"191\u00f79=21, 2\", \"461\u00f77=65, 6\"],\n  [\"998\u00f74=249, 2\", \"640\u00f78=80, 0\"],\n  [\"783\u00f77=111, 6\", \"578\u00f75=115, 3\"],\n  [\"548\u00f72=274, 0\", \"342\u00f74=85, 2\"],\n  [\"728\u00f74=182, 0\", \"420\u00f79=46, 6\"],\n  [\"818\u00f74=204, 2\", \"156\u00f76=26, 0\"],\n  [\"230\u00f76=38, 2\", \"519\u00f74=129, 3\"],\n  [\"290\u00f78=36, 2\", \"120\u00f73=40, 0\"],\n  [\"176\u00f79=19, 5\", \"421\u00f76=70, 1\"],\n  [\"203\u00f72=101, 1\", \"356\u00f74=89, 0\"],\n  [\"563\u00f75=112, 3\", \"569\u00f75=113, 4\"],\n  [\"453\u00f79=50, 3\", \"447\u00f72=223, 1\"],\n  [\"265\u00f79=29, 4\", \"137\u00f78=17, 1\"],\n  [\"948\u00f74=237, 0\", \"789\u00f77=112, 5\"],\n  [\"362\u00f79=40, 2\", \"360\u00f78=45, 0\"],\n  [\"724\u00f75=144, 4\", \"529\u00f78=66, 1\"],\n  [\"828\u00f72=414, 0\", \"635\u00f76=105, 5\"],\n  [\"887\u00f77=126, 5\", \"740\u00f77=105, 5\"],\n  [\"678\u00f75=135, 3\", \"253\u00f77=36, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answers in the table with the new set\n# of problems, per the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"509\u00f74=127, 1\", \"152\u00f73=50, 2\"),\n    @(\"961\u00f79=106, 7\", \"884\u00f77=126, 2\"),\n    @(\"364\u00f79=40, 4\", \"163\u00f76=27, 1\"),\n    @(\"693\u00f78=86, 5\", \"546\u00f79=60, 6\"),\n    @(\"725\u00f72=362, 1\", \"916\u00f75=183, 1\"),\n    @(\"766\u00f72=383, 0\", \"875\u00f79=97, 2\"),\n    @(\"191\u00f79=21, 2\", \"461\u00f77=65, 6\"),\n    @(\"998\u00f74=249, 2\", \"640\u00f78=80, 0\"),\n    @(\"783\u00f77=111, 6\", \"578\u00f75=115, 3\"),\n    @(\"548\u00f72=274, 0\", \"342\u00f74=85, 2\"),\n    @(\"728\u00f74=182, 0\", \"420\u00f79=46, 6\"),\n    @(\"818\u00f74=204, 2\", \"156\u00f76=26, 0\"),\n    @(\"230\u00f76=38, 2\", \"519\u00f74=129, 3\"),\n    @(\"290\u00f78=36, 2\", \"120\u00f73=40, 0\"),\n    @(\"176\u00f79=19, 5\", \"421\u00f76=70, 1\"),\n    @(\"203\u00f72=101, 1\", \"356\u00f74=89, 0\"),\n    @(\"563\u00f75=112, 3\", \"569\u00f75=113, 4\"),\n    @(\"453\u00f79=50, 3\", \"447\u00f72=223, 1\"),\n    @(\"265\u00f79=29, 4\", \"137\u00f78=17, 1\"),\n    @(\"948\u00f74=237, 0\", \"789\u00f77=112, 5\"),\n    @(\"362\u00f79=40, 2\", \"360\u00f78=45, 0\"),\n    @(\"724\u00f75=144, 4\", \"529\u00f78=66, 1\"),\n    @(\"828\u00f72=414, 0\", \"635\u00f76=105, 5\"),\n    @(\"887\u00f77=126, 5\", \"740\u00f77=105, 5\"),\n    @(\"678\u00f75=135, 3\", \"253\u00f77=36, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$true, [ref]$newText, 2)\n}\n"}
